$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "56-27=29"
$t.Cell(1,2).Range.Text = "29+24=53"
$t.Cell(1,3).Range.Text = "32-4=28"
$t.Cell(1,4).Range.Text = "17+77=94"
$t.Cell(1,5).Range.Text = "14+37=51"
$t.Cell(2,1).Range.Text = "17+16=33"
$t.Cell(2,2).Range.Text = "73-54=19"
$t.Cell(2,3).Range.Text = "54-45=9"
$t.Cell(2,4).Range.Text = "44-29=15"
$t.Cell(2,5).Range.Text = "36-8=28"
$t.Cell(3,1).Range.Text = "17+16=33"
$t.Cell(3,2).Range.Text = "64+27=91"
$t.Cell(3,3).Range.Text = "59+35=94"
$t.Cell(3,4).Range.Text = "83-8=75"
$t.Cell(3,5).Range.Text = "46+8=54"
$t.Cell(4,1).Range.Text = "43+9=52"
$t.Cell(4,2).Range.Text = "74-19=55"
$t.Cell(4,3).Range.Text = "64-8=56"
$t.Cell(4,4).Range.Text = "90-59=31"
$t.Cell(4,5).Range.Text = "19+24=43"
$t.Cell(5,1).Range.Text = "30-12=18"
$t.Cell(5,2).Range.Text = "17+68=85"
$t.Cell(5,3).Range.Text = "17+9=26"
$t.Cell(5,4).Range.Text = "76+7=83"
$t.Cell(5,5).Range.Text = "45+36=81"
$t.Cell(6,1).Range.Text = "35+59=94"
$t.Cell(6,2).Range.Text = "39+56=95"
$t.Cell(6,3).Range.Text = "31-29=2"
$t.Cell(6,4).Range.Text = "50-8=42"
$t.Cell(6,5).Range.Text = "81-4=77"
$t.Cell(7,1).Range.Text = "17+69=86"
$t.Cell(7,2).Range.Text = "33+58=91"
$t.Cell(7,3).Range.Text = "13-4=9"
$t.Cell(7,4).Range.Text = "81-23=58"
$t.Cell(7,5).Range.Text = "93-74=19"
$t.Cell(8,1).Range.Text = "54-9=45"
$t.Cell(8,2).Range.Text = "30-12=18"
$t.Cell(8,3).Range.Text = "8+79=87"
$t.Cell(8,4).Range.Text = "18+49=67"
$t.Cell(8,5).Range.Text = "10-7=3"
$t.Cell(9,1).Range.Text = "31-15=16"
$t.Cell(9,2).Range.Text = "62-59=3"
$t.Cell(9,3).Range.Text = "60-54=6"
$t.Cell(9,4).Range.Text = "49+19=68"
$t.Cell(9,5).Range.Text = "70-45=25"
$t.Cell(10,1).Range.Text = "9+6=15"
$t.Cell(10,2).Range.Text = "75-27=48"
$t.Cell(10,3).Range.Text = "28+9=37"
$t.Cell(10,4).Range.Text = "68+5=73"
$t.Cell(10,5).Range.Text = "37+24=61"
$t.Cell(11,1).Range.Text = "66-37=29"
$t.Cell(11,2).Range.Text = "75+7=82"
$t.Cell(11,3).Range.Text = "38+16=54"
$t.Cell(11,4).Range.Text = "55-19=36"
$t.Cell(11,5).Range.Text = "16+19=35"
$t.Cell(12,1).Range.Text = "90-2=88"
$t.Cell(12,2).Range.Text = "93-47=46"
$t.Cell(12,3).Range.Text = "61-55=6"
$t.Cell(12,4).Range.Text = "90-76=14"
$t.Cell(12,5).Range.Text = "80-68=12"
$t.Cell(13,1).Range.Text = "47+36=83"
$t.Cell(13,2).Range.Text = "19+27=46"
$t.Cell(13,3).Range.Text = "41-34=7"
$t.Cell(13,4).Range.Text = "43-18=25"
$t.Cell(13,5).Range.Text = "95-87=8"
$t.Cell(14,1).Range.Text = "33-6=27"
$t.Cell(14,2).Range.Text = "67+26=93"
$t.Cell(14,3).Range.Text = "83-24=59"
$t.Cell(14,4).Range.Text = "16+66=82"
$t.Cell(14,5).Range.Text = "48+18=66"
$t.Cell(15,1).Range.Text = "66-29=37"
$t.Cell(15,2).Range.Text = "97-8=89"
$t.Cell(15,3).Range.Text = "49+13=62"
$t.Cell(15,4).Range.Text = "45-39=6"
$t.Cell(15,5).Range.Text = "90-56=34"
$t.Cell(16,1).Range.Text = "15+36=51"
$t.Cell(16,2).Range.Text = "43-37=6"
$t.Cell(16,3).Range.Text = "10-1=9"
$t.Cell(16,4).Range.Text = "43+19=62"
$t.Cell(16,5).Range.Text = "15+16=31"
$t.Cell(17,1).Range.Text = "22+69=91"
$t.Cell(17,2).Range.Text = "7+29=36"
$t.Cell(17,3).Range.Text = "9+19=28"
$t.Cell(17,4).Range.Text = "6+5=11"
$t.Cell(17,5).Range.Text = "25-6=19"
$t.Cell(18,1).Range.Text = "95-78=17"
$t.Cell(18,2).Range.Text = "29+69=98"
$t.Cell(18,3).Range.Text = "26+19=45"
$t.Cell(18,4).Range.Text = "5+48=53"
$t.Cell(18,5).Range.Text = "31-26=5"
$t.Cell(19,1).Range.Text = "20-2=18"
$t.Cell(19,2).Range.Text = "2+49=51"
$t.Cell(19,3).Range.Text = "55-29=26"
$t.Cell(19,4).Range.Text = "36-29=7"
$t.Cell(19,5).Range.Text = "52-38=14"
$t.Cell(20,1).Range.Text = "72-68=4"
$t.Cell(20,2).Range.Text = "33-27=6"
$t.Cell(20,3).Range.Text = "43+19=62"
$t.Cell(20,4).Range.Text = "73-5=68"
$t.Cell(20,5).Range.Text = "17+59=76"
